$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.814.07"
$ws.Range("D3").Value = "1.567.16"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.03%  "
$r = $ws.Range("D5")
$r.Value = "'206.48"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("E7").Value = "  +0.05%  "
$r = $ws.Range("D8")
$r.Value = "'21.99"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  -1.26%  "
$r = $ws.Range("D11")
$r.Value = "'0.0863"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "1.788.82"
$ws.Range("D13").Value = "1.563.61"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "26.814.24"
$ws.Range("E16").Value = "  -2.40%  "
$r = $ws.Range("D17")
$r.Value = "'61.45"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -3.69%  "
$r = $ws.Range("D18")
$r.Value = "'7.42"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +1.73%  "
$r = $ws.Range("D19")
$r.Value = "'214.64"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  +0.07%  "
$r = $ws.Range("D23")
$r.Value = "'9.28"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -2.93%  "
$r = $ws.Range("D24")
$r.Value = "'1.99"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -0.99%  "
$r = $ws.Range("D25")
$r.Value = "'153.52"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  +0.16%  "
$r = $ws.Range("D27")
$r.Value = "'14.98"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "1.396.34"
$ws.Range("E34").Value = "  -1.83%  "
$r = $ws.Range("D35")
$r.Value = "'1.53"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("E39").Value = "  -2.96%  "
$r = $ws.Range("D40")
$r.Value = "'0.816"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("E41").Value = "  +0.06%  "
$r = $ws.Range("D42")
$r.Value = "'0.991"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +1.01%  "
$r = $ws.Range("D43")
$r.Value = "'1.80"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("E45").Value = "  +0.79%  "
$r = $ws.Range("D46")
$r.Value = "'63.33"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").Value = "1.702.60"
$ws.Range("E47").Value = "  -0.04%  "
$r = $ws.Range("D48")
$r.Value = "'86.06"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("D49").Value = "0.0₇0986"
$ws.Range("E49").Value = "  -1.46%  "
$r = $ws.Range("D50")
$r.Value = "'0.0953"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("E51").Value = "  -0.91%  "
